$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 799.82355
$ws.Range("I19").Value = 634.7
$ws.Range("J19").Value = 1035.7142
$ws.Range("K19").Value = 634.7
$ws.Range("L19").Value = 1035.7142
$ws.Range("M19").Value = -459.7
$ws.Range("N19").Value = -1385.7142

$ws.Range("H51").Value = 4116.567
$ws.Range("J51").Value = 4411.294
$ws.Range("L51").Value = 4411.294
$ws.Range("N51").Value = -5379.294

$ws.Range("H70").Value = 4246
$ws.Range("J70").Value = 3993
$ws.Range("L70").Value = 11979
$ws.Range("N70").Value = -12519

$ws.Range("H73").Value = 4246
$ws.Range("J73").Value = 3993
$ws.Range("L73").Value = 11979
$ws.Range("N73").Value = -13851

$ws.Range("H74").Value = 75071
$ws.Range("I74").Value = 86916.586
$ws.Range("J74").Value = 3997.5
$ws.Range("K74").Value = 86916.586
$ws.Range("L74").Value = 3997.5
$ws.Range("M74").Value = -85980.586
$ws.Range("N74").Value = -5869.5

$ws.Range("H77").Value = 75071
$ws.Range("I77").Value = 86916.586
$ws.Range("J77").Value = 3997.5
$ws.Range("K77").Value = 434582.93
$ws.Range("L77").Value = 19987.5
$ws.Range("M77").Value = -429902.93
$ws.Range("N77").Value = -29347.5

$ws.Range("H98").Value = 1311.7407
$ws.Range("I98").Value = 1348.3334
$ws.Range("K98").Value = 1348.3334
$ws.Range("M98").Value = 149.6666

$ws.Range("H111").Value = 1366.7
$ws.Range("I111").Value = 1292.7778
$ws.Range("J111").Value = 2032
$ws.Range("K111").Value = 3878.3334
$ws.Range("L111").Value = 6096
$ws.Range("M111").Value = -811.3334000000004
$ws.Range("N111").Value = -12230

$ws.Range("H112").Value = 2918.3125
$ws.Range("I112").Value = 2352.4443
$ws.Range("K112").Value = 7057.3329
$ws.Range("M112").Value = -5949.3329

$ws.Range("H113").Value = 2299.3333
$ws.Range("I113").Value = 2299.3333
$ws.Range("K113").Value = 2299.3333
$ws.Range("M113").Value = 954.6667000000002

$ws.Range("H116").Value = 36459.375
$ws.Range("J116").Value = 40260
$ws.Range("L116").Value = 40260
$ws.Range("N116").Value = -47144

$ws.Range("H122").Value = 1311.7407
$ws.Range("I122").Value = 1348.3334
$ws.Range("K122").Value = 4045.0002
$ws.Range("M122").Value = -1595.0002

$ws.Range("H138").Value = 3651.9
$ws.Range("J138").Value = 3554.4814
$ws.Range("L138").Value = 10663.4442
$ws.Range("N138").Value = -20943.4442

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3575944.8
$ws.Range("I61").Value = 4559.923
$ws.Range("J61").Value = 50003948
$ws.Range("K61").Value = 4559.923
$ws.Range("L61").Value = 50003948
$ws.Range("M61").Value = -4347.923
$ws.Range("N61").Value = -50004372

$ws.Range("H63").Value = 5062.375
$ws.Range("I63").Value = 5874.75
$ws.Range("K63").Value = 5874.75
$ws.Range("M63").Value = -5188.75

$ws.Range("H66").Value = 5062.375
$ws.Range("I66").Value = 5874.75
$ws.Range("K66").Value = 29373.75
$ws.Range("M66").Value = -25941.75

$ws.Range("H103").Value = 37776
$ws.Range("J103").Value = 37776
$ws.Range("L103").Value = 37776
$ws.Range("N103").Value = -40120

$ws.Range("H122").Value = 2260
$ws.Range("I122").Value = 1765
$ws.Range("K122").Value = 5295
$ws.Range("M122").Value = -2845

$ws.Range("H132").Value = 964050.4399999999
$ws.Range("I132").Value = 1089396.1
$ws.Range("J132").Value = 3066.6667
$ws.Range("K132").Value = 3268188.3
$ws.Range("L132").Value = 9200.000100000001
$ws.Range("M132").Value = -3265658.3
$ws.Range("N132").Value = -14260.0001

$ws.Range("H136").Value = 3575944.8
$ws.Range("I136").Value = 4559.923
$ws.Range("J136").Value = 50003948
$ws.Range("K136").Value = 13679.769
$ws.Range("L136").Value = 150011844
$ws.Range("M136").Value = -11129.769
$ws.Range("N136").Value = -150016944

$ws.Range("H140").Value = 73800
$ws.Range("J140").Value = 73800
$ws.Range("L140").Value = 73800
$ws.Range("N140").Value = -84160

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 8
$ws.Range("I14").Value = 8
$ws.Range("K14").Value = 8
$ws.Range("M14").Value = 164

$ws.Range("H86").Value = 3003
$ws.Range("I86").Value = 3003
$ws.Range("K86").Value = 3003
$ws.Range("M86").Value = -1880

$ws.Range("H89").Value = 3003
$ws.Range("I89").Value = 3003
$ws.Range("K89").Value = 15015
$ws.Range("M89").Value = -9399

$ws.Range("H97").Value = 34284
$ws.Range("I97").Value = 21438.75
$ws.Range("J97").Value = 59974.5
$ws.Range("K97").Value = 21438.75
$ws.Range("L97").Value = 59974.5
$ws.Range("M97").Value = -20447.75
$ws.Range("N97").Value = -61956.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1157.4667
$ws.Range("I22").Value = 236.7
$ws.Range("K22").Value = 236.7
$ws.Range("M22").Value = 113.3

$ws.Range("H105").Value = 7912.1875
$ws.Range("I105").Value = 10382.546
$ws.Range("J105").Value = 2477.4
$ws.Range("K105").Value = 10382.546
$ws.Range("L105").Value = 2477.4
$ws.Range("M105").Value = -8635.546
$ws.Range("N105").Value = -5971.4

$ws.Range("H122").Value = 14325
$ws.Range("I122").Value = 2944.625
$ws.Range("K122").Value = 8833.875
$ws.Range("M122").Value = -6383.875

$ws.Range("H132").Value = 3289.2
$ws.Range("I132").Value = 3199.2632
$ws.Range("K132").Value = 9597.7896
$ws.Range("M132").Value = -7067.7896

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 111.3
$ws.Range("I2").Value = 111.3
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 667.8
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -554.8
$ws.Range("N2").ClearContents()

$ws.Range("H38").Value = 97.59999999999999
$ws.Range("I38").Value = 72.166664
$ws.Range("J38").Value = 135.75
$ws.Range("K38").Value = 216.499992
$ws.Range("L38").Value = 407.25
$ws.Range("M38").Value = 130.500008
$ws.Range("N38").Value = -1101.25

$ws.Range("H92").Value = 357.66666
$ws.Range("J92").Value = 338
$ws.Range("L92").Value = 1014
$ws.Range("N92").Value = -3510

$ws.Range("H131").Value = 6622.5454
$ws.Range("J131").Value = 8643.75
$ws.Range("L131").Value = 25931.25
$ws.Range("N131").Value = -36011.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3240.3333
$ws.Range("I113").Value = 2785.7334
$ws.Range("J113").Value = 3998
$ws.Range("K113").Value = 2785.7334
$ws.Range("L113").Value = 3998
$ws.Range("M113").Value = -615.7334000000001
$ws.Range("N113").Value = -8338

$ws.Range("H126").Value = 2439.8333
$ws.Range("I126").Value = 2356.25
$ws.Range("K126").Value = 7068.75
$ws.Range("M126").Value = -4598.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1134.6923
$ws.Range("I16").Value = 337.67648
$ws.Range("K16").Value = 337.67648
$ws.Range("M16").Value = -167.67648

$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 15550
$ws.Range("I62").Value = 15000
$ws.Range("J62").Value = 15733.333
$ws.Range("K62").Value = 15000
$ws.Range("L62").Value = 15733.333
$ws.Range("M62").Value = -14376
$ws.Range("N62").Value = -16981.333

$ws.Range("H65").Value = 15550
$ws.Range("I65").Value = 15000
$ws.Range("J65").Value = 15733.333
$ws.Range("K65").Value = 75000
$ws.Range("L65").Value = 78666.66500000001
$ws.Range("M65").Value = -71880
$ws.Range("N65").Value = -84906.66500000001

$ws.Range("H81").Value = 1888.7142
$ws.Range("I81").Value = 1370.1666
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 2740.3332
$ws.Range("L81").Value = 10000
$ws.Range("M81").Value = -1679.3332
$ws.Range("N81").Value = -12122

$ws.Range("H84").Value = 1888.7142
$ws.Range("I84").Value = 1370.1666
$ws.Range("J84").Value = 5000
$ws.Range("K84").Value = 13701.666
$ws.Range("L84").Value = 50000
$ws.Range("M84").Value = -8397.666000000001
$ws.Range("N84").Value = -60608

$ws.Range("H135").Value = 105238.336
$ws.Range("J135").Value = 105238.336
$ws.Range("L135").Value = 105238.336
$ws.Range("N135").Value = -115378.336
